# Auto update Excel log
# Appends newly-logged sensor events to four of the SeniorConnect_MasterLog
# worksheets: ALERTS, PIR, Humidity and Proximity.
#
# Every data cell in this workbook is stored as literal text (dates like
# "2026-01-30" and times like "16:44:56" are NOT real date/time values), so
# each cell is forced to Text format ("@") before the value is written --
# otherwise Excel would silently reinterpret strings such as "2026-01-30"
# as a date serial or "87.7%" as a number.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $rowNum, $colA, $colB, $colC, $colD, $colE, $colF)

    $vals = @($colA, $colB, $colC, $colD, $colE, $colF)
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $vals[$c - 1]
    }
}

# ---------------------------------------------------------------------
# ALERTS sheet -- one new MINIMAL-motion alert row
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")

Add-LogRow $wsAlerts 2 "2026-01-30" "16:44:56" "16:00" "Bathroom" "MINIMAL" "MINIMAL ALERT: Bathroom occupied, no motion > 20s."

# ---------------------------------------------------------------------
# PIR sheet -- 14 new "No Motion" / "Inactive" Bathroom readings (rows 22-35)
# ---------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")

$pirTimes = @("16:44:00","16:44:00","16:44:03","16:44:09","16:44:13","16:44:19","16:44:23","16:44:29","16:44:34","16:44:39","16:44:44","16:44:49","16:44:54","16:44:59")

for ($i = 0; $i -lt $pirTimes.Count; $i++) {
    $rowNum = 22 + $i
    Add-LogRow $wsPir $rowNum "2026-01-30" $pirTimes[$i] "16:00" "Bathroom" "No Motion" "Inactive"
}

# ---------------------------------------------------------------------
# Humidity sheet -- 8 new "87.7%" / "Active" Bathroom readings (rows 14-21)
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityTimes = @("16:44:00","16:44:04","16:44:19","16:44:24","16:44:29","16:44:39","16:44:44","16:44:59")

for ($i = 0; $i -lt $humidityTimes.Count; $i++) {
    $rowNum = 14 + $i
    Add-LogRow $wsHumidity $rowNum "2026-01-30" $humidityTimes[$i] "16:00" "Bathroom" "87.7%" "Active"
}

# ---------------------------------------------------------------------
# Proximity sheet -- 5 new Bathroom Door ENTER/EXIT events (rows 14-18)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

$proxTimes  = @("16:44:03","16:44:07","16:44:20","16:44:32","16:44:53")
$proxEvents = @("EXIT","ENTER","EXIT","ENTER","EXIT")
$proxStatus = @("User EXITED Bathroom","User ENTERED Bathroom","User EXITED Bathroom","User ENTERED Bathroom","User EXITED Bathroom")

for ($i = 0; $i -lt $proxTimes.Count; $i++) {
    $rowNum = 14 + $i
    Add-LogRow $wsProximity $rowNum "2026-01-30" $proxTimes[$i] "16:00" "Bathroom Door" $proxEvents[$i] $proxStatus[$i]
}
